$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.084.68"
$ws.Range("E2").Value = "  +3.19%  "
$ws.Range("D3").Value = "1.674.24"
$ws.Range("E3").Value = "  +1.64%  "
$ws.Range("E4").Value = "  -0.74%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9990"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3646"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.31"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3242"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.146"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07218"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9984"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.078"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.95%  "
$ws.Range("D15").Value = "1.672.29"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.661"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001052"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06531"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9990"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.887"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.87%  "
$ws.Range("D24").Value = "25.066.34"
$ws.Range("E24").Value = "  +2.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.430"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.388"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "149.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.18%  "
$ws.Range("D29").Value = "1.856.33"
$ws.Range("E29").Value = "  +1.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.192"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.073"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.783"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08425"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.653"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.164"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.234"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.95%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06053"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.80%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02220"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2077"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.292"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9994"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5957"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.62"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.827"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5703"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.956"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07145"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.181"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.85%  "
